# "just creat the new LUR input"
# Copy the HOA/COA (B/C) values and polygon id (E) values that already
# live in rows 2-11 of the "polygon" sheet into a second, side-by-side
# table anchored at D15:E25 (with polygon id moving to column C), giving
# the new LUR input block.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("200m")
$ws2 = $wb.Worksheets.Item("polygon")

# Header row for the new block: D15 = "HOA", E15 = "COA"
$ws2.Range("D15").Value = "HOA"
$ws2.Range("E15").Value = "COA"

# Source data already present on the sheet: B2:C11 (HOA/COA) and E2:E11
# (polygon id), now duplicated into C16:E25.
for ($i = 0; $i -lt 10; $i++) {
    $srcRow = 2 + $i
    $dstRow = 16 + $i

    $polygonId = $ws2.Cells.Item($srcRow, 5).Value2   # column E
    $hoa = $ws2.Cells.Item($srcRow, 2).Value2          # column B
    $coa = $ws2.Cells.Item($srcRow, 3).Value2          # column C

    $ws2.Cells.Item($dstRow, 3).Value = $polygonId    # column C
    $ws2.Cells.Item($dstRow, 4).Value = $hoa           # column D
    $ws2.Cells.Item($dstRow, 5).Value = $coa           # column E
}

# Restore the view state captured in the saved workbook: "200m" is the
# tab that ends up selected/active, with the "polygon" sheet scrolled so
# its new block is in view.
$ws2.Range("E19").Select()
$ws1.Range("E51").Select()
$ws1.Activate()
